$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B label updates and Column C count updates
$ws.Range("B2").Value = "<up>"
$ws.Range("C2").Value = 31

$ws.Range("B3").Value = "<shult>"

$ws.Range("B4").Value = "<been>"
$ws.Range("C4").Value = 31

$ws.Range("B5").Value = "<fima>"
$ws.Range("C5").Value = 40

$ws.Range("C6").Value = 37

$ws.Range("C7").Value = 26

$ws.Range("B8").Value = "<lipa>"
$ws.Range("C8").Value = 36

$ws.Range("B9").Value = "<echi>"
$ws.Range("C9").Value = 42

$ws.Range("B10").Value = "<cot>"
$ws.Range("C10").Value = 34

$ws.Range("B11").Value = "<firold>"
$ws.Range("C11").Value = 35

$ws.Range("B12").Value = "<had>"
$ws.Range("C12").Value = 32

$ws.Range("C13").Value = 36

$ws.Range("C14").Value = 33

$ws.Range("C15").Value = 40

$ws.Range("B16").Value = "<ech>"
$ws.Range("C16").Value = 34

$ws.Range("C17").Value = 33

$ws.Range("C18").Value = 35
